$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells remain text (avoid Excel auto-converting numbers/dates)
$targetCells = @("D2","E2","D3","E3","E4","D5","E5","D6","E6","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","D20","E20","D21","E21","D22","E22","D23","E23","D24","E24","B25","C25","D25","E25","B26","C26","D26","E26","D27","E27","D28","E28","B29","C29","D29","E29","B30","C30","D30","E30","D31","E31","D32","E32","D33","E33","D34","E34","D35","E35","D36","E36","D37","E37","D38","E38","D39","E39","B40","C40","D40","E40","B41","C41","D41","E41","B42","C42","D42","E42","B43","C43","D43","E43","D44","E44","D45","E45","D46","E46","D47","E47","B48","C48","D48","E48","B49","C49","D49","E49","B50","C50","D50","E50","D51","E51")
foreach ($addr in $targetCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values from the diff
$ws.Range("D2").Value = '61.292.67'
$ws.Range("E2").Value = '  -5.28%  '
$ws.Range("D3").Value = '2.950.91'
$ws.Range("E3").Value = '  -7.15%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '536.54'
$ws.Range("E5").Value = '  -6.08%  '
$ws.Range("D6").Value = '151.98'
$ws.Range("E6").Value = '  -7.86%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("D8").Value = '0.558'
$ws.Range("E8").Value = '  -5.10%  '
$ws.Range("D9").Value = '2.954.79'
$ws.Range("E9").Value = '  -6.83%  '
$ws.Range("D10").Value = '0.111'
$ws.Range("E10").Value = '  -5.71%  '
$ws.Range("D11").Value = '6.10'
$ws.Range("E11").Value = '  -8.60%  '
$ws.Range("D12").Value = '0.363'
$ws.Range("E12").Value = '  -6.43%  '
$ws.Range("D13").Value = '3.467.64'
$ws.Range("E13").Value = '  -7.05%  '
$ws.Range("D14").Value = '0.124'
$ws.Range("E14").Value = '  -3.30%  '
$ws.Range("D15").Value = '61.312.21'
$ws.Range("E15").Value = '  -5.39%  '
$ws.Range("D16").Value = '23.56'
$ws.Range("E16").Value = '  -7.69%  '
$ws.Range("D17").Value = '2.952.85'
$ws.Range("E17").Value = '  -6.79%  '
$ws.Range("D18").Value = '0.0000145'
$ws.Range("E18").Value = '  -7.42%  '
$ws.Range("D19").Value = '5.10'
$ws.Range("E19").Value = '  -3.90%  '
$ws.Range("D20").Value = '377.72'
$ws.Range("E20").Value = '  -7.89%  '
$ws.Range("D21").Value = '11.88'
$ws.Range("E21").Value = '  -7.18%  '
$ws.Range("D22").Value = '6.60'
$ws.Range("E22").Value = '  -7.58%  '
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("D24").Value = '64.72'
$ws.Range("E24").Value = '  -6.22%  '
$ws.Range("B25").Value = 'Polygon'
$ws.Range("C25").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D25").Value = '0.464'
$ws.Range("E25").Value = '  -4.94%  '
$ws.Range("B26").Value = 'WrappedeETH'
$ws.Range("C26").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D26").Value = '3.066.57'
$ws.Range("E26").Value = '  -7.67%  '
$ws.Range("D27").Value = '0.185'
$ws.Range("E27").Value = '  -7.91%  '
$ws.Range("D28").Value = '0.998'
$ws.Range("E28").Value = '  -0.30%  '
$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").Value = '0.0₃0918'
$ws.Range("E29").Value = '  -11.76%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").Value = '8.25'
$ws.Range("E30").Value = '  -7.40%  '
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  +0.00%  '
$ws.Range("D32").Value = '1.70'
$ws.Range("E32").Value = '  -6.68%  '
$ws.Range("D33").Value = '20.15'
$ws.Range("E33").Value = '  -5.52%  '
$ws.Range("D34").Value = '156.17'
$ws.Range("E34").Value = '  -0.76%  '
$ws.Range("D35").Value = '5.94'
$ws.Range("E35").Value = '  -7.11%  '
$ws.Range("D36").Value = '4.56'
$ws.Range("E36").Value = '  -7.84%  '
$ws.Range("D37").Value = '1.05'
$ws.Range("E37").Value = '  -7.73%  '
$ws.Range("D38").Value = '1.27'
$ws.Range("E38").Value = '  -6.66%  '
$ws.Range("D39").Value = '1.53'
$ws.Range("E39").Value = '  -10.54%  '
$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D40").Value = '3.88'
$ws.Range("E40").Value = '  -5.41%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '2.401.05'
$ws.Range("E41").Value = '  -11.47%  '
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").Value = '22.00'
$ws.Range("E42").Value = '  -8.72%  '
$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").Value = '36.64'
$ws.Range("E43").Value = '  -5.73%  '
$ws.Range("D44").Value = '0.658'
$ws.Range("E44").Value = '  -5.98%  '
$ws.Range("D45").Value = '0.0588'
$ws.Range("E45").Value = '  -6.14%  '
$ws.Range("D46").Value = '0.997'
$ws.Range("E46").Value = '  -0.20%  '
$ws.Range("D47").Value = '0.0242'
$ws.Range("E47").Value = '  -7.02%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = '4.90'
$ws.Range("E48").Value = '  -10.05%  '
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").Value = '0.0947'
$ws.Range("E49").Value = '  -3.90%  '
$ws.Range("B50").Value = 'WhiteBITCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D50").Value = '10.50'
$ws.Range("E50").Value = '  +0.06%  '
$ws.Range("D51").Value = '19.49'
$ws.Range("E51").Value = '  -9.53%  '
